$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the fund / portfolio company / amount values for row 2
$ws.Range("A2").Value = "Demo Fund 2"
$ws.Range("B2").Value = "TSTF2 Port Co 3"
$ws.Range("D2").Value = 10000000

# Update the selected cell/range on the sheet
$ws.Range("C2").Select()
